$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row 1 cells: "_old" suffix -> "_FV2304", "_new" suffix -> "_FV2310"
$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the header range into an Excel Table (ListObject) so headers carry the
# new names into xl/tables/table1.xml as well.
$range = $ws.Range("A1:U81")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (pane split after row 1).
$ws.Activate()
[void]$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
